$d = $word.ActiveDocument

# Locate the paragraph that still contains the grammar-flagged phrase
# "in depth" (wrapped in <w:proofErr w:type="gramStart"/> / gramEnd markers).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "*in depth*") {
        $target = $para.Range
        break
    }
}

if ($target -eq $null) {
    throw "Could not find paragraph containing 'in depth'"
}

# Rebuild the paragraph's OOXML: turn "in depth" into the single word
# "in-depth" and drop the two now-unneeded <w:proofErr/> grammar-check
# markers, leaving every other run untouched.
$rsquo = [char]0x2019
$paraXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
    'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' +
    'w14:paraId="1F58997C" w14:textId="5F8D06C5" w:rsidR="000247F8" w:rsidRDefault="000247F8">' +
    '<w:r><w:t xml:space="preserve">For a more detailed comparison of our ideal jobs, </w:t></w:r>' +
    '<w:r w:rsidR="00B43BA2"><w:t xml:space="preserve">the industry data section of the Fighting Mongoose' + $rsquo + 's homepage provides more </w:t></w:r>' +
    '<w:r><w:t>in-depth</w:t></w:r>' +
    '<w:r w:rsidR="00B43BA2"><w:t xml:space="preserve"> analysis using data from Burning Glass </w:t></w:r>' +
    '<w:r w:rsidR="00572063"><w:t>Technologies</w:t></w:r>' +
    '<w:r w:rsidR="00B43BA2"><w:t>.</w:t></w:r>' +
    '</w:p>'

$target.InsertXML($paraXml) | Out-Null
